# Add meaningful alt text to the QR code picture in the instructions
# template. The document has a single inline picture (the "Survey QR
# Code" graphic) sitting right under the cover-sheet merge field; it was
# previously auto-tagged by Word's "Generate Alt Text" feature with a
# generic/auto-generated description. Replace it with a short, useful
# description and make sure the run is flagged "do not spell/grammar
# check" the way Word marks runs that contain only a picture.

$d = $word.ActiveDocument

$qr = $d.InlineShapes.Item(1)

# Updates both wp:docPr/@descr and pic:cNvPr/@descr for the picture.
$qr.AlternativeText = "Survey QR Code"

# Mirrors Word's usual <w:rPr><w:noProof/></w:rPr> on a run that only
# contains a picture/drawing.
$qr.Range.NoProofing = $true
